# Auto-generated: update live market-price derived columns (H-N)
# on the Leve-profit sheets, matching a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 11476.889
$ws.Range("I34").Value = 11476.889
$ws.Range("K34").Value = 11476.889
$ws.Range("M34").Value = -11273.889
$ws.Range("H36").Value = 11476.889
$ws.Range("I36").Value = 11476.889
$ws.Range("K36").Value = 11476.889
$ws.Range("M36").Value = -10761.889
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H51").Value = 2959.6775
$ws.Range("I51").Value = 2575
$ws.Range("J51").Value = 3659.0908
$ws.Range("K51").Value = 2575
$ws.Range("L51").Value = 3659.0908
$ws.Range("M51").Value = -2091
$ws.Range("N51").Value = -4627.0908
$ws.Range("H64").Value = 8281.1875
$ws.Range("I64").Value = 5727.1816
$ws.Range("J64").Value = 13900
$ws.Range("K64").Value = 5727.1816
$ws.Range("L64").Value = 13900
$ws.Range("M64").Value = -5479.1816
$ws.Range("N64").Value = -14396
$ws.Range("H67").Value = 8281.1875
$ws.Range("I67").Value = 5727.1816
$ws.Range("J67").Value = 13900
$ws.Range("K67").Value = 5727.1816
$ws.Range("L67").Value = 13900
$ws.Range("M67").Value = -4869.1816
$ws.Range("N67").Value = -15616
$ws.Range("H98").Value = 3694.2964
$ws.Range("I98").Value = 1593.8422
$ws.Range("J98").Value = 8682.875
$ws.Range("K98").Value = 1593.8422
$ws.Range("L98").Value = 8682.875
$ws.Range("M98").Value = -95.84220000000005
$ws.Range("N98").Value = -11678.875
$ws.Range("H99").Value = 201.25
$ws.Range("I99").Value = 201.25
$ws.Range("K99").Value = 603.75
$ws.Range("M99").Value = 894.25
$ws.Range("H112").Value = 53235.81
$ws.Range("J112").Value = 61797.945
$ws.Range("L112").Value = 185393.835
$ws.Range("N112").Value = -187609.835
$ws.Range("H122").Value = 3694.2964
$ws.Range("I122").Value = 1593.8422
$ws.Range("J122").Value = 8682.875
$ws.Range("K122").Value = 4781.5266
$ws.Range("L122").Value = 26048.625
$ws.Range("M122").Value = -2331.5266
$ws.Range("N122").Value = -30948.625
$ws.Range("H125").Value = 3171.7646
$ws.Range("I125").Value = 3120.625
$ws.Range("J125").Value = 3217.2222
$ws.Range("K125").Value = 28085.625
$ws.Range("L125").Value = 28954.9998
$ws.Range("M125").Value = -25625.625
$ws.Range("N125").Value = -33874.99980000001
$ws.Range("H135").Value = 1166.12
$ws.Range("I135").Value = 1110.5416
$ws.Range("K135").Value = 9994.874400000001
$ws.Range("M135").Value = -7459.874400000001
$ws.Range("H138").Value = 3974.26
$ws.Range("I138").Value = 1780.6842
$ws.Range("J138").Value = 4488.8022
$ws.Range("K138").Value = 5342.0526
$ws.Range("L138").Value = 13466.4066
$ws.Range("M138").Value = -202.0526
$ws.Range("N138").Value = -23746.4066
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 16489.908
$ws.Range("I45").Value = 24484.572
$ws.Range("J45").Value = 2499.25
$ws.Range("K45").Value = 24484.572
$ws.Range("L45").Value = 2499.25
$ws.Range("M45").Value = -24107.572
$ws.Range("N45").Value = -3253.25
$ws.Range("H61").Value = 189089.67
$ws.Range("I61").Value = 2240.7058
$ws.Range("J61").Value = 894963.5600000001
$ws.Range("K61").Value = 2240.7058
$ws.Range("L61").Value = 894963.5600000001
$ws.Range("M61").Value = -2028.7058
$ws.Range("N61").Value = -895387.5600000001
$ws.Range("H97").Value = 1806.5
$ws.Range("J97").Value = 862.3333
$ws.Range("L97").Value = 862.3333
$ws.Range("N97").Value = -1854.3333
$ws.Range("H132").Value = 2239.0667
$ws.Range("I132").Value = 1896.7805
$ws.Range("J132").Value = 5747.5
$ws.Range("K132").Value = 5690.3415
$ws.Range("L132").Value = 17242.5
$ws.Range("M132").Value = -3160.3415
$ws.Range("N132").Value = -22302.5
$ws.Range("H136").Value = 189089.67
$ws.Range("I136").Value = 2240.7058
$ws.Range("J136").Value = 894963.5600000001
$ws.Range("K136").Value = 6722.117400000001
$ws.Range("L136").Value = 2684890.68
$ws.Range("M136").Value = -4172.117400000001
$ws.Range("N136").Value = -2689990.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3476.9375
$ws.Range("I20").Value = 3392.3333
$ws.Range("K20").Value = 3392.3333
$ws.Range("M20").Value = -3145.3333
$ws.Range("H99").Value = 4615.32
$ws.Range("I99").Value = 3077.3572
$ws.Range("K99").Value = 3077.3572
$ws.Range("M99").Value = -1579.3572
$ws.Range("H134").Value = 1709.7858
$ws.Range("I134").Value = 1709.7858
$ws.Range("K134").Value = 5129.357400000001
$ws.Range("M134").Value = -2594.357400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1634.8148
$ws.Range("I105").Value = 1531.1364
$ws.Range("K105").Value = 1531.1364
$ws.Range("M105").Value = 215.8635999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 438.7143
$ws.Range("I75").Value = 352.75
$ws.Range("K75").Value = 1058.25
$ws.Range("M75").Value = -60.25
$ws.Range("H78").Value = 438.7143
$ws.Range("I78").Value = 352.75
$ws.Range("K78").Value = 3174.75
$ws.Range("M78").Value = 1817.25
$ws.Range("H86").Value = 1033.3334
$ws.Range("I86").Value = 800
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 2400
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -1214
$ws.Range("N86").Value = -6872
$ws.Range("H89").Value = 1033.3334
$ws.Range("I89").Value = 800
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 7200
$ws.Range("L89").Value = 13500
$ws.Range("M89").Value = -1272
$ws.Range("N89").Value = -25356

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2693.4333
$ws.Range("J22").Value = 3474.75
$ws.Range("L22").Value = 3474.75
$ws.Range("N22").Value = -4064.75
$ws.Range("H27").Value = 2693.4333
$ws.Range("J27").Value = 3474.75
$ws.Range("L27").Value = 3474.75
$ws.Range("N27").Value = -3688.75
$ws.Range("H46").Value = 3424.762
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 3713.75
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 3713.75
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -4089.75
$ws.Range("H68").Value = 2941.3076
$ws.Range("J68").Value = 2514.2222
$ws.Range("L68").Value = 2514.2222
$ws.Range("N68").Value = -4012.2222
$ws.Range("H71").Value = 2941.3076
$ws.Range("J71").Value = 2514.2222
$ws.Range("L71").Value = 12571.111
$ws.Range("N71").Value = -20059.111
$ws.Range("H82").Value = 12654.9
$ws.Range("J82").Value = 3833.3333
$ws.Range("L82").Value = 3833.3333
$ws.Range("N82").Value = -4555.3333
$ws.Range("H85").Value = 12654.9
$ws.Range("J85").Value = 3833.3333
$ws.Range("L85").Value = 3833.3333
$ws.Range("N85").Value = -6329.3333
$ws.Range("H122").Value = 7321.9395
$ws.Range("I122").Value = 7213.1904
$ws.Range("J122").Value = 7512.25
$ws.Range("K122").Value = 21639.5712
$ws.Range("L122").Value = 22536.75
$ws.Range("M122").Value = -19189.5712
$ws.Range("N122").Value = -27436.75
$ws.Range("H132").Value = 4833.5713
$ws.Range("I132").Value = 4333.3335
$ws.Range("K132").Value = 13000.0005
$ws.Range("M132").Value = -10470.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 40555
$ws.Range("J27").Value = 40555
$ws.Range("L27").Value = 40555
$ws.Range("N27").Value = -40693
$ws.Range("H46").Value = 97666.664
$ws.Range("J46").Value = 97666.664
$ws.Range("L46").Value = 97666.664
$ws.Range("N46").Value = -98128.664
$ws.Range("H62").Value = 19993
$ws.Range("I62").Value = 14989
$ws.Range("K62").Value = 14989
$ws.Range("M62").Value = -14365
$ws.Range("H65").Value = 19993
$ws.Range("I65").Value = 14989
$ws.Range("K65").Value = 74945
$ws.Range("M65").Value = -71825
$ws.Range("H107").Value = 1339.5
$ws.Range("J107").Value = 1823.75
$ws.Range("L107").Value = 5471.25
$ws.Range("N107").Value = -9311.25
$ws.Range("H122").Value = 1942.3715
$ws.Range("I122").Value = 1860.742
$ws.Range("K122").Value = 5582.226
$ws.Range("M122").Value = -3132.226
$ws.Range("H130").Value = 18036.334
$ws.Range("J130").Value = 18036.334
$ws.Range("L130").Value = 18036.334
$ws.Range("N130").Value = -28076.334
$ws.Range("H134").Value = 97666.664
$ws.Range("J134").Value = 97666.664
$ws.Range("L134").Value = 292999.992
$ws.Range("N134").Value = -298069.992
$ws.Range("H136").Value = 2121.1936
$ws.Range("I136").Value = 1848.4642
$ws.Range("J136").Value = 4666.6665
$ws.Range("K136").Value = 5545.392599999999
$ws.Range("L136").Value = 13999.9995
$ws.Range("M136").Value = -2995.392599999999
$ws.Range("N136").Value = -19099.9995
